# Apply cryptos list update (prices & volume %) scraped Sat Jan 13 21:10:06 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.886.38'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.02%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.567.43'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.28'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.72'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.31%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.94%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.549'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.70'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.00%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0811'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.86%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.72'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.31%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.577.54'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.40%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.885'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.38'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '42.901.27'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.91%  '

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0996'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.17%  '

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.85'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.63'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.01'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.92%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '254.51'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.23%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.30%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.23%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '28.74'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.75%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.26'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.98'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.67%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.11'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -5.56%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.03'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.45%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '154.96'

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'LidoDAOToken'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.40'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -5.08%  '

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'WEMIXToken'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.75'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.26%  '

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.15'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.95%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0802'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.41%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.18'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +8.48%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.113'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.18%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.93%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.04'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.42'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.74%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0310'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.87'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +25.64%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.068.68'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.18'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.24'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '76.58'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +10.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '106.45'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.98%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.823.62'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.191'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.27%  '
